$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9; this pushes the existing row 9
# (and everything below it) down to row 10, copying formatting along.
$ws.Rows("9:9").Insert()

# The newly inserted row 9 is currently a duplicate of the old row 9
# (now at row 10). Copy the old row's values/formats into the new row,
# then overwrite the fields that changed for the new weekly entry.
$ws.Rows("10:10").Copy()
$ws.Rows("9:9").PasteSpecial()

# Update the cells that changed in the new (latest) entry.
$ws.Range("D9").Value = 45015
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 24000
$ws.Range("M9").Value = 24000
$ws.Range("P9").Value = 1333
